$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) Overview sheet: status text changed from "Ready for handoff" to
#    "Handed back: in sync with en-US" for both language columns / both rows.
# ---------------------------------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("B2").Value = "Handed back: in sync with en-US"
$overview.Range("C2").Value = "Handed back: in sync with en-US"
$overview.Range("B3").Value = "Handed back: in sync with en-US"
$overview.Range("C3").Value = "Handed back: in sync with en-US"

# ---------------------------------------------------------------------------
# 2) zh-cn sheet: populate "Latest Target File" (F) and "Latest Handback
#    File" (G) columns for both data rows, with hyperlinks that mirror the
#    existing "Latest Handoff File" values/links, and record the actual
#    handback datetime in column H (replacing the placeholder date).
# ---------------------------------------------------------------------------
$zhcn = $wb.Worksheets.Item("zh-cn")

$zhcn.Hyperlinks.Add($zhcn.Range("F2"), "https://github.com/OpenLocalizationTest/oltest/blob/122c0d6b61fb7a31e543c3eac7b1112b3ec5d2fd/e2e/7fcf7547-11b3-4378-a2df-fb30585a65b1.md", "", "", "7fcf7547-11b3-4378-a2df-fb30585a65b1.md")
$zhcn.Hyperlinks.Add($zhcn.Range("G2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/10e81ff892d9d619853fbb19ae8f185e6ab4dec1/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/high/7fcf7547-11b3-4378-a2df-fb30585a65b1.17da252eb469da3f3770090ccb1c98c491995c6f.zh-cn.xlf", "", "", "7fcf7547-11b3-4378-a2df-fb30585a65b1.17da252eb469da3f3770090ccb1c98c491995c6f.zh-cn.xlf")
$zhcn.Hyperlinks.Add($zhcn.Range("F3"), "https://github.com/OpenLocalizationTest/oltest/blob/122c0d6b61fb7a31e543c3eac7b1112b3ec5d2fd/e2e/97137af2-7f70-4394-a854-94ec4a75fcf0.md", "", "", "97137af2-7f70-4394-a854-94ec4a75fcf0.md")
$zhcn.Hyperlinks.Add($zhcn.Range("G3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/10e81ff892d9d619853fbb19ae8f185e6ab4dec1/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/high/97137af2-7f70-4394-a854-94ec4a75fcf0.36b55d94ee2c051f2c206e670b25b760ac547d44.zh-cn.xlf", "", "", "97137af2-7f70-4394-a854-94ec4a75fcf0.36b55d94ee2c051f2c206e670b25b760ac547d44.zh-cn.xlf")

$zhcn.Range("F2").Style = "HyperLink"
$zhcn.Range("G2").Style = "HyperLink"
$zhcn.Range("F3").Style = "HyperLink"
$zhcn.Range("G3").Style = "HyperLink"

$zhcn.Range("H2").Value = "2016-03-17 18:10:56"
$zhcn.Range("H3").Value = "2016-03-17 18:10:56"

# ---------------------------------------------------------------------------
# 3) de-de sheet: same shape of change as zh-cn, but with its own handback
#    datetime value.
# ---------------------------------------------------------------------------
$dede = $wb.Worksheets.Item("de-de")

$dede.Hyperlinks.Add($dede.Range("F2"), "https://github.com/OpenLocalizationTest/oltest/blob/122c0d6b61fb7a31e543c3eac7b1112b3ec5d2fd/e2e/7fcf7547-11b3-4378-a2df-fb30585a65b1.md", "", "", "7fcf7547-11b3-4378-a2df-fb30585a65b1.md")
$dede.Hyperlinks.Add($dede.Range("G2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/5c646358667ea831390a72a3a29bff6081c2b584/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/high/7fcf7547-11b3-4378-a2df-fb30585a65b1.17da252eb469da3f3770090ccb1c98c491995c6f.de-de.xlf", "", "", "7fcf7547-11b3-4378-a2df-fb30585a65b1.17da252eb469da3f3770090ccb1c98c491995c6f.de-de.xlf")
$dede.Hyperlinks.Add($dede.Range("F3"), "https://github.com/OpenLocalizationTest/oltest/blob/122c0d6b61fb7a31e543c3eac7b1112b3ec5d2fd/e2e/97137af2-7f70-4394-a854-94ec4a75fcf0.md", "", "", "97137af2-7f70-4394-a854-94ec4a75fcf0.md")
$dede.Hyperlinks.Add($dede.Range("G3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/5c646358667ea831390a72a3a29bff6081c2b584/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/high/97137af2-7f70-4394-a854-94ec4a75fcf0.36b55d94ee2c051f2c206e670b25b760ac547d44.de-de.xlf", "", "", "97137af2-7f70-4394-a854-94ec4a75fcf0.36b55d94ee2c051f2c206e670b25b760ac547d44.de-de.xlf")

$dede.Range("F2").Style = "HyperLink"
$dede.Range("G2").Style = "HyperLink"
$dede.Range("F3").Style = "HyperLink"
$dede.Range("G3").Style = "HyperLink"

$dede.Range("H2").Value = "2016-03-17 18:11:04"
$dede.Range("H3").Value = "2016-03-17 18:11:04"

$wb.Save()
